# priming update (look at why compareKeys isn't working)
#
# The author was investigating the "correct_key" column (J) by
# highlighting the "A" (apple-associated) answers in yellow and the "L"
# (horse-associated) answers in green for a couple of representative
# blocks, and duplicated the practice block's block_number values into a
# new helper column K. An AutoFilter was also dropped onto the
# correct_key column while debugging.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New helper column K: copy of block_number (A) for the practice rows ---
# Using Range.Copy(destination) preserves both the value (reusing the
# existing "practice" shared string) and the existing cell style (s="1"),
# without registering any new style/font.
$ws.Range("A2:A10").Copy($ws.Range("K2:K10"))

# --- Highlight certain correct_key (J) cells ---
# Yellow fill (standard "Yellow", FFFF00) on the "A" rows being inspected.
$ws.Range("J2").Interior.Color = 65535
$ws.Range("J6").Interior.Color = 65535
$ws.Range("J10:J16").Interior.Color = 65535
$ws.Range("J35:J40").Interior.Color = 65535

# Green fill (standard "Green", 92D050) on the "L" rows being inspected.
$ws.Range("J23:J28").Interior.Color = 5296274
$ws.Range("J47:J52").Interior.Color = 5296274

# --- Turn on an AutoFilter on the correct_key column ---
$ws.Range("J1:J59").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$J`$1:`$J`$59")
$filterName.Visible = $false

# --- Update the active selection to where the author was last looking ---
$ws.Range("I17").Select()
